# Set detector offset and gain to zero for cryo biases.
#
# The "Biases" sheet has one row per PAR value (column A) with bias
# settings across columns B:BO. Rows 34-41 (PAR 99-106) previously held a
# flat default of 1000 across every column; zero them out instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Biases")

$ws.Range("B34:BO41").Value = 0

# Reflect the author's new scroll position / active cell in the sheet view.
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1

$ws.Range("A34").Select()
